$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.711.25'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  -4.68%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.649.80'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  -6.54%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  +1.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.16'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -2.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3638'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -5.55%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3293'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  -9.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.75'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  -8.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.134'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -8.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07106'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -7.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.051'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  -6.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.76'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -9.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.670'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -6.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.647.90'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -6.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001064'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -8.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06597'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  -3.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.003'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '79.67'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -8.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.37'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -7.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.045'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -7.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.25'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -4.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.730.00'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -4.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.412'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -0.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.563'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -14.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '148.94'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -4.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.33'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -7.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '127.90'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -5.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.835.42'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -6.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.132'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -6.46%  '
$ws.Range('E32').Value = '  -4.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.133'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -14.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.723'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -5.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08473'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -2.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.73'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -9.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.224'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -8.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06228'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -8.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02290'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -8.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2111'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -5.71%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.228'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -5.98%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.339'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -10.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6114'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -7.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.16'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -7.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.762'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -4.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5779'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -9.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.992'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -8.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.52'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -7.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07062'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -5.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '75.13'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -7.06%  '
